# Adapt column header formatting to respective input file names:
# "<header>_old" -> "<header>_FV2410" and "<header>_new" -> "<header>_FV2504",
# then expose the data range as an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
    "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410", "Bedingung_FV2410", "diff",
    "Segmentname_FV2504", "Segmentgruppe_FV2504", "Segment_FV2504", "Datenelement_FV2504",
    "Segment ID_FV2504", "Code_FV2504", "Qualifier_FV2504", "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504", "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table ("Table1") covering A1:U64.
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header row renamed, Table1 created, panes frozen."
